# Apply "correction of python code with new sources" edits to the sources sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sources")

# ---------------------------------------------------------------------------
# 1. Column B (Source) -- most WSJ quote URLs now point at the
#    historical-prices sub-page; GOLD gets a trailing slash; the Brent
#    futures link and the three global-rates.com rate links are unchanged.
# ---------------------------------------------------------------------------
$sources = @{
    2  = "https://www.wsj.com/market-data/quotes/index/SPX/historical-prices"
    3  = "https://www.wsj.com/market-data/quotes/index/DJIA/historical-prices"
    4  = "https://www.wsj.com/market-data/quotes/index/UK/UKX/historical-prices"
    5  = "https://www.wsj.com/market-data/quotes/index/DX/DAX/historical-prices"
    6  = "https://www.wsj.com/market-data/quotes/index/FR/PX1/historical-prices"
    7  = "https://www.wsj.com/market-data/quotes/index/XX/XSTX/SX5E/historical-prices"
    8  = "https://www.wsj.com/market-data/quotes/index/CN/SHCOMP/historical-prices"
    9  = "https://www.wsj.com/market-data/quotes/index/JP/XTKS/NIK/historical-prices"
    10 = "https://www.wsj.com/market-data/quotes/index/XX/990100/historical-prices"
    11 = "https://www.wsj.com/market-data/quotes/index/XX/891800/historical-prices"
    12 = "https://www.wsj.com/market-data/quotes/index/VIX/historical-prices"
    13 = "https://www.wsj.com/market-data/quotes/futures/UK/BRENT%20CRUDE"
    14 = "https://www.wsj.com/market-data/quotes/futures/GOLD/"
    15 = "https://www.wsj.com/market-data/quotes/fx/EURUSD/historical-prices"
    16 = "https://www.wsj.com/market-data/quotes/fx/USDJPY/historical-prices"
    17 = "https://www.wsj.com/market-data/quotes/fx/GBPUSD/historical-prices"
    18 = "https://www.wsj.com/market-data/quotes/fx/USDCHF/historical-prices"
    19 = "https://www.global-rates.com/interest-rates/eonia/eonia.aspx"
    20 = "https://www.global-rates.com/interest-rates/libor/american-dollar/usd-libor-interest-rate-3-months.aspx"
    21 = "https://www.global-rates.com/interest-rates/euribor/euribor-interest-3-months.aspx"
    22 = "https://www.wsj.com/market-data/quotes/bond/BX/TMUBMUSD10Y/historical-prices"
    23 = "https://www.wsj.com/market-data/quotes/bond/BX/TMBMKDE-10Y/historical-prices"
    24 = "https://www.wsj.com/market-data/quotes/bond/BX/TMBMKFR-10Y/historical-prices"
    25 = "https://www.wsj.com/market-data/quotes/bond/BX/TMBMKJP-10Y/historical-prices"
    26 = "https://www.wsj.com/market-data/quotes/bond/BX/TMBMKIT-10Y/historical-prices"
}

foreach ($row in $sources.Keys) {
    $ws.Cells.Item($row, 2).Value2 = $sources[$row]
}

# ---------------------------------------------------------------------------
# 2. Column C (Init) -- the eight rate rows (19-26) were re-expressed as
#    whole percentage points instead of decimal fractions (x100).
# ---------------------------------------------------------------------------
$rateRows = 19..26
foreach ($row in $rateRows) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value2 = $cell.Value2 * 100
}

# ---------------------------------------------------------------------------
# 3. Hyperlinks -- rebuild the whole collection so every link's target
#    matches the new Source text (this also drops the stale cached
#    "display" text that Excel had kept for the Brent Crude link).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()
foreach ($row in $sources.Keys) {
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 2), $sources[$row]) | Out-Null
}

# ---------------------------------------------------------------------------
# 4. Cosmetic workbook state: page setup + current selection.
# ---------------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("A2").Select()
